# -----------------------------------------------------------------------
# Commit: "rename E03 & added thumnails"
#   1. Slide 2 ("01. Introduction") textbox is widened to full slide
#      width, centered, word-wrap enabled, and re-titled to
#      "01. Motivation".
#   2. Five new slides are appended (clones of slide 3's "02. Structure"
#      section-divider layout), titled:
#        03. Getting Started
#        04. stdCallback
#        05. stdLambda
#        06. stdArray
#        07. stdEnumerator
# -----------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Slide 2: "01. Introduction" -> "01. Motivation" ----------------
$s2 = $p.Slides.Item(2)
$tb = $s2.Shapes.Item(3)

$tb.Name = "TextBox 7"
$tb.Left = 0
$tb.Width = 960
$tb.TextFrame.WordWrap = $true
$tb.TextFrame.TextRange.Text = "01. Motivation"
$tb.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- 2. Append 5 new "section divider" slides, cloned from slide 3 -----
#    Each entry is (firstRunText, secondRunTextOrNull)
$titles = @(
    @("03. Getting Started", $null),
    @("04. ", "stdCallback"),
    @("05. ", "stdLambda"),
    @("06. ", "stdArray"),
    @("07. ", "stdEnumerator")
)

$src = $p.Slides.Item(3)

foreach ($pair in $titles) {
    $firstText = $pair[0]
    $secondText = $pair[1]

    $dup = $src.Duplicate()
    $newSlide = $dup.Item(1)
    $ntb = $newSlide.Shapes.Item(3)
    $tr = $ntb.TextFrame.TextRange

    # Clearing the (possibly multi-run) inherited text first keeps the
    # subsequent assignment from being mis-split across the old runs.
    $tr.Text = ""
    $tr.Text = $firstText
    if ($secondText) {
        $tr.InsertAfter($secondText) | Out-Null
    }

    # chain off the slide we just built, so the next duplicate lands
    # immediately after it (keeping presentation order == title order)
    $src = $newSlide
}
